# Apply the updates described by the diff to the "Metadata" sheet.
# Because the "Title" and "Description" values are also duplicated verbatim
# on the "Elements" sheet (Short / Definition columns for the root Extension
# row), and the diff only edits the shared-string text itself (not the
# cell-to-string mappings), those matching cells must be updated too so the
# text stays shared/deduplicated exactly like before.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# Version: 1.0.1 -> 0.0.0
$meta.Range("B3").Value = "0.0.0"

# Title: "Extension of Patient Ethnicity" -> "Ethnicity"
$meta.Range("B5").Value = "Ethnicity"

# Date: 2023-11-21T19:08:35-03:00 -> 2024-01-16T20:08:55-03:00
$meta.Range("B8").Value = "2024-01-16T20:08:55-03:00"

# Description: drop the trailing ICHOM sentence
$meta.Range("B12").Value = "Extension to capture the patient's ethnicity which represents their cultural background or heritage."

# The "Elements" sheet repeats the Title text in the "Short" column (L2)
# and the Description text in the "Definition" column (M2) for the root
# Extension row - keep them in sync with the same shared-string values.
$elements.Range("L2").Value = "Ethnicity"
$elements.Range("M2").Value = "Extension to capture the patient's ethnicity which represents their cultural background or heritage."
